# Auto-generated: applies per-cell market/profit data refresh
# as produced by the scheduled Sheets runner (see commit message).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(15, 8).Value = 2200.2964  # H15: 2491.1875 -> 2200.2964
$ws.Cells.Item(15, 9).Value = 2200.2964  # I15: 2491.1875 -> 2200.2964
$ws.Cells.Item(15, 11).Value = 6600.889200000001  # K15: 7473.5625 -> 6600.889200000001
$ws.Cells.Item(15, 13).Value = -6431.889200000001  # M15: -7304.5625 -> -6431.889200000001
$ws.Cells.Item(74, 8).Value = 5657.88  # H74: 4718 -> 5657.88
$ws.Cells.Item(74, 9).Value = 4851.9585  # I74: 4718 -> 4851.9585
$ws.Cells.Item(74, 10).Value = 25000  # J74: 0 -> 25000
$ws.Cells.Item(74, 11).Value = 4851.9585  # K74: 4718 -> 4851.9585
$ws.Cells.Item(74, 12).Value = 25000  # L74: 0 -> 25000
$ws.Cells.Item(74, 13).Value = -3915.9585  # M74: -3782 -> -3915.9585
$ws.Cells.Item(74, 14).Value = -26872  # N74: None -> -26872
$ws.Cells.Item(77, 8).Value = 5657.88  # H77: 4718 -> 5657.88
$ws.Cells.Item(77, 9).Value = 4851.9585  # I77: 4718 -> 4851.9585
$ws.Cells.Item(77, 10).Value = 25000  # J77: 0 -> 25000
$ws.Cells.Item(77, 11).Value = 24259.7925  # K77: 23590 -> 24259.7925
$ws.Cells.Item(77, 12).Value = 125000  # L77: 0 -> 125000
$ws.Cells.Item(77, 13).Value = -19579.7925  # M77: -18910 -> -19579.7925
$ws.Cells.Item(77, 14).Value = -134360  # N77: None -> -134360
$ws.Cells.Item(80, 8).Value = 21067.312  # H80: 22385.2 -> 21067.312
$ws.Cells.Item(80, 9).Value = 11174.917  # I80: 12072.728 -> 11174.917
$ws.Cells.Item(80, 11).Value = 33524.751  # K80: 36218.18399999999 -> 33524.751
$ws.Cells.Item(80, 13).Value = -32526.751  # M80: -35220.18399999999 -> -32526.751
$ws.Cells.Item(83, 8).Value = 21067.312  # H83: 22385.2 -> 21067.312
$ws.Cells.Item(83, 9).Value = 11174.917  # I83: 12072.728 -> 11174.917
$ws.Cells.Item(83, 11).Value = 100574.253  # K83: 108654.552 -> 100574.253
$ws.Cells.Item(83, 13).Value = -95582.253  # M83: -103662.552 -> -95582.253
$ws.Cells.Item(103, 8).Value = 949.75  # H103: 933 -> 949.75
$ws.Cells.Item(103, 10).Value = 949.75  # J103: 933 -> 949.75
$ws.Cells.Item(103, 12).Value = 2849.25  # L103: 2799 -> 2849.25
$ws.Cells.Item(103, 14).Value = -4021.25  # N103: -3971 -> -4021.25
$ws.Cells.Item(112, 8).Value = 1950.2941  # H112: 1953.0555 -> 1950.2941
$ws.Cells.Item(112, 10).Value = 1953.75  # J112: 1956.4706 -> 1953.75
$ws.Cells.Item(112, 12).Value = 5861.25  # L112: 5869.4118 -> 5861.25
$ws.Cells.Item(112, 14).Value = -8077.25  # N112: -8085.4118 -> -8077.25
$ws.Cells.Item(134, 8).Value = 65687.06  # H134: 63943.555 -> 65687.06
$ws.Cells.Item(134, 10).Value = 65687.06  # J134: 63943.555 -> 65687.06
$ws.Cells.Item(134, 12).Value = 65687.06  # L134: 63943.555 -> 65687.06
$ws.Cells.Item(134, 14).Value = -75827.06  # N134: -74083.55499999999 -> -75827.06
$ws.Cells.Item(137, 8).Value = 27634.312  # H137: 17746.96 -> 27634.312
$ws.Cells.Item(137, 9).Value = 1751  # I137: 1700.3334 -> 1751
$ws.Cells.Item(137, 10).Value = 31331.928  # J137: 19840 -> 31331.928
$ws.Cells.Item(137, 11).Value = 5253  # K137: 5101.0002 -> 5253
$ws.Cells.Item(137, 12).Value = 93995.784  # L137: 59520 -> 93995.784
$ws.Cells.Item(137, 13).Value = -2703  # M137: -2551.0002 -> -2703
$ws.Cells.Item(137, 14).Value = -99095.784  # N137: -64620 -> -99095.784
$ws.Cells.Item(138, 8).Value = 5130.0786  # H138: 5111.981 -> 5130.0786
$ws.Cells.Item(138, 10).Value = 6653.606  # J138: 6581.1177 -> 6653.606
$ws.Cells.Item(138, 12).Value = 19960.818  # L138: 19743.3531 -> 19960.818
$ws.Cells.Item(138, 14).Value = -30240.818  # N138: -30023.3531 -> -30240.818

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 2129.5833  # H32: 1071.974 -> 2129.5833
$ws.Cells.Item(32, 9).Value = 2266.2188  # I32: 1132.884 -> 2266.2188
$ws.Cells.Item(32, 10).Value = 1036.5  # J32: 546.625 -> 1036.5
$ws.Cells.Item(32, 11).Value = 2266.2188  # K32: 1132.884 -> 2266.2188
$ws.Cells.Item(32, 12).Value = 1036.5  # L32: 546.625 -> 1036.5
$ws.Cells.Item(32, 13).Value = -1979.2188  # M32: -845.884 -> -1979.2188
$ws.Cells.Item(32, 14).Value = -1610.5  # N32: -1120.625 -> -1610.5
$ws.Cells.Item(37, 8).Value = 43206.6  # H37: 39008.25 -> 43206.6
$ws.Cells.Item(37, 9).Value = 39008.25  # I37: 32011 -> 39008.25
$ws.Cells.Item(37, 11).Value = 39008.25  # K37: 32011 -> 39008.25
$ws.Cells.Item(37, 13).Value = -38735.25  # M37: -31738 -> -38735.25
$ws.Cells.Item(44, 8).Value = 50000  # H44: 40000 -> 50000
$ws.Cells.Item(44, 10).Value = 50000  # J44: 40000 -> 50000
$ws.Cells.Item(44, 12).Value = 50000  # L44: 40000 -> 50000
$ws.Cells.Item(44, 14).Value = -50976  # N44: -40976 -> -50976
$ws.Cells.Item(45, 8).Value = 54937.58  # H45: 52260.2 -> 54937.58
$ws.Cells.Item(45, 9).Value = 64425.875  # I45: 60717.883 -> 64425.875
$ws.Cells.Item(45, 11).Value = 64425.875  # K45: 60717.883 -> 64425.875
$ws.Cells.Item(45, 13).Value = -64048.875  # M45: -60340.883 -> -64048.875
$ws.Cells.Item(55, 8).Value = 30499.5  # H55: 34666.332 -> 30499.5
$ws.Cells.Item(61, 8).Value = 6282.8335  # H61: 4487.857 -> 6282.8335
$ws.Cells.Item(61, 9).Value = 6400  # I61: 4030.2727 -> 6400
$ws.Cells.Item(61, 11).Value = 6400  # K61: 4030.2727 -> 6400
$ws.Cells.Item(61, 13).Value = -6188  # M61: -3818.2727 -> -6188
$ws.Cells.Item(80, 8).Value = 79397.5  # H80: 79400 -> 79397.5
$ws.Cells.Item(80, 10).Value = 79397.5  # J80: 79400 -> 79397.5
$ws.Cells.Item(80, 12).Value = 79397.5  # L80: 79400 -> 79397.5
$ws.Cells.Item(80, 14).Value = -81393.5  # N80: -81396 -> -81393.5
$ws.Cells.Item(83, 8).Value = 79397.5  # H83: 79400 -> 79397.5
$ws.Cells.Item(83, 10).Value = 79397.5  # J83: 79400 -> 79397.5
$ws.Cells.Item(83, 12).Value = 238192.5  # L83: 238200 -> 238192.5
$ws.Cells.Item(83, 14).Value = -248176.5  # N83: -248184 -> -248176.5
$ws.Cells.Item(132, 8).Value = 457903.75  # H132: 422611.12 -> 457903.75
$ws.Cells.Item(132, 9).Value = 981990.5600000001  # I132: 771828.1 -> 981990.5600000001
$ws.Cells.Item(132, 10).Value = 14445.692  # J132: 15191.25 -> 14445.692
$ws.Cells.Item(132, 11).Value = 2945971.68  # K132: 2315484.3 -> 2945971.68
$ws.Cells.Item(132, 12).Value = 43337.076  # L132: 45573.75 -> 43337.076
$ws.Cells.Item(132, 13).Value = -2943441.68  # M132: -2312954.3 -> -2943441.68
$ws.Cells.Item(132, 14).Value = -48397.076  # N132: -50633.75 -> -48397.076
$ws.Cells.Item(136, 8).Value = 6282.8335  # H136: 4487.857 -> 6282.8335
$ws.Cells.Item(136, 9).Value = 6400  # I136: 4030.2727 -> 6400
$ws.Cells.Item(136, 11).Value = 19200  # K136: 12090.8181 -> 19200
$ws.Cells.Item(136, 13).Value = -16650  # M136: -9540.8181 -> -16650

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 2159.3333  # H20: 2159 -> 2159.3333
$ws.Cells.Item(20, 9).Value = 1399.6364  # I20: 1399.1818 -> 1399.6364
$ws.Cells.Item(20, 11).Value = 1399.6364  # K20: 1399.1818 -> 1399.6364
$ws.Cells.Item(20, 13).Value = -1152.6364  # M20: -1152.1818 -> -1152.6364
$ws.Cells.Item(86, 8).Value = 2666.6667  # H86: 3000 -> 2666.6667
$ws.Cells.Item(89, 8).Value = 2666.6667  # H89: 3000 -> 2666.6667
$ws.Cells.Item(99, 8).Value = 3333.3333  # H99: 3500 -> 3333.3333
$ws.Cells.Item(99, 9).Value = 3333.3333  # I99: 3500 -> 3333.3333
$ws.Cells.Item(99, 11).Value = 3333.3333  # K99: 3500 -> 3333.3333
$ws.Cells.Item(99, 13).Value = -1835.3333  # M99: -2002 -> -1835.3333
$ws.Cells.Item(107, 8).Value = 803.6875  # H107: 833.93335 -> 803.6875
$ws.Cells.Item(107, 9).Value = 623.93335  # I107: 643.5 -> 623.93335
$ws.Cells.Item(107, 11).Value = 623.93335  # K107: 643.5 -> 623.93335
$ws.Cells.Item(107, 13).Value = 1296.06665  # M107: 1276.5 -> 1296.06665

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 483536.4  # H31: 414776.44 -> 483536.4
$ws.Cells.Item(31, 9).Value = 737750.6  # I31: 707052.7 -> 737750.6
$ws.Cells.Item(31, 10).Value = 33772.77  # J31: 25074.723 -> 33772.77
$ws.Cells.Item(31, 11).Value = 737750.6  # K31: 707052.7 -> 737750.6
$ws.Cells.Item(31, 12).Value = 33772.77  # L31: 25074.723 -> 33772.77
$ws.Cells.Item(31, 13).Value = -737455.6  # M31: -706757.7 -> -737455.6
$ws.Cells.Item(31, 14).Value = -34362.77  # N31: -25664.723 -> -34362.77
$ws.Cells.Item(34, 8).Value = 483536.4  # H34: 414776.44 -> 483536.4
$ws.Cells.Item(34, 9).Value = 737750.6  # I34: 707052.7 -> 737750.6
$ws.Cells.Item(34, 10).Value = 33772.77  # J34: 25074.723 -> 33772.77
$ws.Cells.Item(34, 11).Value = 737750.6  # K34: 707052.7 -> 737750.6
$ws.Cells.Item(34, 12).Value = 33772.77  # L34: 25074.723 -> 33772.77
$ws.Cells.Item(34, 13).Value = -737548.6  # M34: -706850.7 -> -737548.6
$ws.Cells.Item(34, 14).Value = -34176.77  # N34: -25478.723 -> -34176.77
$ws.Cells.Item(62, 8).Value = 4833.3335  # H62: 6250 -> 4833.3335
$ws.Cells.Item(62, 9).Value = 2000  # I62: 0 -> 2000
$ws.Cells.Item(62, 11).Value = 2000  # K62: 0 -> 2000
$ws.Cells.Item(62, 13).Value = -1376  # M62: None -> -1376
$ws.Cells.Item(65, 8).Value = 4833.3335  # H65: 6250 -> 4833.3335
$ws.Cells.Item(65, 9).Value = 2000  # I65: 0 -> 2000
$ws.Cells.Item(65, 11).Value = 10000  # K65: 0 -> 10000
$ws.Cells.Item(65, 13).Value = -6880  # M65: None -> -6880
$ws.Cells.Item(132, 8).Value = 3176.8928  # H132: 3313.6538 -> 3176.8928
$ws.Cells.Item(132, 9).Value = 1997.9  # I132: 2064.4443 -> 1997.9
$ws.Cells.Item(132, 11).Value = 5993.700000000001  # K132: 6193.3329 -> 5993.700000000001
$ws.Cells.Item(132, 13).Value = -3463.700000000001  # M132: -3663.3329 -> -3463.700000000001
$ws.Cells.Item(134, 8).Value = 275894.06  # H134: 309199.12 -> 275894.06
$ws.Cells.Item(134, 9).Value = 2288.25  # I134: 2481.75 -> 2288.25
$ws.Cells.Item(134, 11).Value = 6864.75  # K134: 7445.25 -> 6864.75
$ws.Cells.Item(134, 13).Value = -4329.75  # M134: -4910.25 -> -4329.75

# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(2, 8).Value = 2618.4375  # H2: 157.91667 -> 2618.4375
$ws.Cells.Item(2, 9).Value = 5049.1665  # I2: 98.333336 -> 5049.1665
$ws.Cells.Item(2, 10).Value = 1160  # J2: 177.77777 -> 1160
$ws.Cells.Item(2, 11).Value = 30294.999  # K2: 590.000016 -> 30294.999
$ws.Cells.Item(2, 12).Value = 6960  # L2: 1066.66662 -> 6960
$ws.Cells.Item(2, 13).Value = -30181.999  # M2: -477.000016 -> -30181.999
$ws.Cells.Item(2, 14).Value = -7186  # N2: -1292.66662 -> -7186
$ws.Cells.Item(3, 8).Value = 8499.75  # H3: 10199.6 -> 8499.75
$ws.Cells.Item(3, 9).Value = 6999.6665  # I3: 9499.5 -> 6999.6665
$ws.Cells.Item(3, 11).Value = 20998.9995  # K3: 28498.5 -> 20998.9995
$ws.Cells.Item(3, 13).Value = -20886.9995  # M3: -28386.5 -> -20886.9995
$ws.Cells.Item(80, 8).Value = 0  # H80: 1590 -> 0
$ws.Cells.Item(80, 9).Value = 0  # I80: 1590 -> 0
$ws.Cells.Item(80, 11).Value = 0  # K80: 4770 -> 0
$ws.Cells.Item(80, 13).Value = ""  # M80 (was -3834) -> cleared
$ws.Cells.Item(83, 8).Value = 0  # H83: 1590 -> 0
$ws.Cells.Item(83, 9).Value = 0  # I83: 1590 -> 0
$ws.Cells.Item(83, 11).Value = 0  # K83: 14310 -> 0
$ws.Cells.Item(83, 13).Value = ""  # M83 (was -9630) -> cleared
$ws.Cells.Item(122, 8).Value = 12690882  # H122: 13666950 -> 12690882
$ws.Cells.Item(122, 10).Value = 1556.7  # J122: 1507.4445 -> 1556.7
$ws.Cells.Item(122, 12).Value = 14010.3  # L122: 13567.0005 -> 14010.3
$ws.Cells.Item(122, 14).Value = -18910.3  # N122: -18467.0005 -> -18910.3
$ws.Cells.Item(129, 8).Value = 875  # H129: 1102 -> 875
$ws.Cells.Item(129, 9).Value = 245.4  # I129: 276.33334 -> 245.4
$ws.Cells.Item(129, 10).Value = 1924.3334  # J129: 1455.8572 -> 1924.3334
$ws.Cells.Item(129, 11).Value = 736.2  # K129: 829.0000200000001 -> 736.2
$ws.Cells.Item(129, 12).Value = 5773.0002  # L129: 4367.571599999999 -> 5773.0002
$ws.Cells.Item(129, 13).Value = 4263.8  # M129: 4170.99998 -> 4263.8
$ws.Cells.Item(129, 14).Value = -15773.0002  # N129: -14367.5716 -> -15773.0002
$ws.Cells.Item(134, 8).Value = 1755.8  # H134: 1315.4286 -> 1755.8
$ws.Cells.Item(134, 9).Value = 1755.8  # I134: 1315.4286 -> 1755.8
$ws.Cells.Item(134, 11).Value = 5267.4  # K134: 3946.2858 -> 5267.4
$ws.Cells.Item(134, 13).Value = -197.3999999999996  # M134: 1123.7142 -> -197.3999999999996
$ws.Cells.Item(141, 8).Value = 4100  # H141: 4333.3335 -> 4100
$ws.Cells.Item(141, 10).Value = 5250  # J141: 5333.3335 -> 5250
$ws.Cells.Item(141, 12).Value = 15750  # L141: 16000.0005 -> 15750
$ws.Cells.Item(141, 14).Value = -26110  # N141: -26360.0005 -> -26110

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 122.521736  # H2: 127.77273 -> 122.521736
$ws.Cells.Item(2, 9).Value = 127.59091  # I2: 133.33333 -> 127.59091
$ws.Cells.Item(2, 11).Value = 127.59091  # K2: 133.33333 -> 127.59091
$ws.Cells.Item(2, 13).Value = -14.59090999999999  # M2: -20.33332999999999 -> -14.59090999999999
$ws.Cells.Item(4, 8).Value = 0  # H4: 20000000 -> 0
$ws.Cells.Item(4, 9).Value = 0  # I4: 20000000 -> 0
$ws.Cells.Item(4, 11).Value = 0  # K4: 20000000 -> 0
$ws.Cells.Item(4, 13).Value = ""  # M4 (was -19999888) -> cleared
$ws.Cells.Item(12, 8).Value = 15000  # H12: 15000000 -> 15000
$ws.Cells.Item(12, 10).Value = 15000  # J12: 15000000 -> 15000
$ws.Cells.Item(12, 12).Value = 15000  # L12: 15000000 -> 15000
$ws.Cells.Item(12, 14).Value = -15280  # N12: -15000280 -> -15280
$ws.Cells.Item(58, 8).Value = 31333.334  # H58: 32000 -> 31333.334
$ws.Cells.Item(58, 9).Value = 30500  # I58: 31500 -> 30500
$ws.Cells.Item(58, 11).Value = 30500  # K58: 31500 -> 30500
$ws.Cells.Item(58, 13).Value = -30223  # M58: -31223 -> -30223
$ws.Cells.Item(97, 8).Value = 508.75  # H97: 1200 -> 508.75
$ws.Cells.Item(97, 9).Value = 445.9091  # I97: 0 -> 445.9091
$ws.Cells.Item(97, 11).Value = 445.9091  # K97: 0 -> 445.9091
$ws.Cells.Item(97, 13).Value = 50.09089999999998  # M97: None -> 50.09089999999998
$ws.Cells.Item(123, 8).Value = 42399.2  # H123: 40285.145 -> 42399.2
$ws.Cells.Item(123, 10).Value = 42399.2  # J123: 40285.145 -> 42399.2
$ws.Cells.Item(123, 12).Value = 42399.2  # L123: 40285.145 -> 42399.2
$ws.Cells.Item(123, 14).Value = -47299.2  # N123: -45185.145 -> -47299.2
$ws.Cells.Item(132, 8).Value = 51196.555  # H132: 53139.117 -> 51196.555
$ws.Cells.Item(132, 9).Value = 15512.833  # I132: 16157.305 -> 15512.833
$ws.Cells.Item(132, 11).Value = 46538.499  # K132: 48471.915 -> 46538.499
$ws.Cells.Item(132, 13).Value = -44008.499  # M132: -45941.915 -> -44008.499
$ws.Cells.Item(134, 8).Value = 87163  # H134: 0 -> 87163
$ws.Cells.Item(134, 10).Value = 87163  # J134: 0 -> 87163
$ws.Cells.Item(134, 12).Value = 261489  # L134: 0 -> 261489
$ws.Cells.Item(134, 14).Value = -266559  # N134: None -> -266559

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(40, 8).Value = 409306  # H40: 464389 -> 409306
$ws.Cells.Item(40, 9).Value = 464143.88  # I40: 510028.4 -> 464143.88
$ws.Cells.Item(40, 10).Value = 7161.6665  # J40: 7995 -> 7161.6665
$ws.Cells.Item(40, 11).Value = 464143.88  # K40: 510028.4 -> 464143.88
$ws.Cells.Item(40, 12).Value = 7161.6665  # L40: 7995 -> 7161.6665
$ws.Cells.Item(40, 13).Value = -464007.88  # M40: -509892.4 -> -464007.88
$ws.Cells.Item(40, 14).Value = -7433.6665  # N40: -8267 -> -7433.6665
$ws.Cells.Item(93, 8).Value = 3824.8333  # H93: 4120.5 -> 3824.8333
$ws.Cells.Item(93, 9).Value = 3824.8333  # I93: 4120.5 -> 3824.8333
$ws.Cells.Item(93, 11).Value = 3824.8333  # K93: 4120.5 -> 3824.8333
$ws.Cells.Item(93, 13).Value = -2576.8333  # M93: -2872.5 -> -2576.8333
$ws.Cells.Item(132, 8).Value = 4148.381  # H132: 4070.5117 -> 4148.381
$ws.Cells.Item(132, 9).Value = 2991.9092  # I132: 2927.4412 -> 2991.9092
$ws.Cells.Item(132, 11).Value = 8975.7276  # K132: 8782.3236 -> 8975.7276
$ws.Cells.Item(132, 13).Value = -6445.7276  # M132: -6252.3236 -> -6445.7276

# --- Sheet: WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(81, 8).Value = 2965.4075  # H81: 2902.4814 -> 2965.4075
$ws.Cells.Item(81, 9).Value = 2305.25  # I81: 2135.6155 -> 2305.25
$ws.Cells.Item(81, 10).Value = 3493.5334  # J81: 3614.5715 -> 3493.5334
$ws.Cells.Item(81, 11).Value = 4610.5  # K81: 4271.231 -> 4610.5
$ws.Cells.Item(81, 12).Value = 6987.0668  # L81: 7229.143 -> 6987.0668
$ws.Cells.Item(81, 13).Value = -3549.5  # M81: -3210.231 -> -3549.5
$ws.Cells.Item(81, 14).Value = -9109.066800000001  # N81: -9351.143 -> -9109.066800000001
$ws.Cells.Item(84, 8).Value = 2965.4075  # H84: 2902.4814 -> 2965.4075
$ws.Cells.Item(84, 9).Value = 2305.25  # I84: 2135.6155 -> 2305.25
$ws.Cells.Item(84, 10).Value = 3493.5334  # J84: 3614.5715 -> 3493.5334
$ws.Cells.Item(84, 11).Value = 23052.5  # K84: 21356.155 -> 23052.5
$ws.Cells.Item(84, 12).Value = 34935.334  # L84: 36145.715 -> 34935.334
$ws.Cells.Item(84, 13).Value = -17748.5  # M84: -16052.155 -> -17748.5
$ws.Cells.Item(84, 14).Value = -45543.334  # N84: -46753.715 -> -45543.334
$ws.Cells.Item(130, 8).Value = 83163.336  # H130: 84740 -> 83163.336
$ws.Cells.Item(130, 10).Value = 83163.336  # J130: 84740 -> 83163.336
$ws.Cells.Item(130, 12).Value = 83163.336  # L130: 84740 -> 83163.336
$ws.Cells.Item(130, 14).Value = -93203.336  # N130: -94780 -> -93203.336
$ws.Cells.Item(132, 8).Value = 25791.775  # H132: 17523.85 -> 25791.775
$ws.Cells.Item(132, 9).Value = 1917.7667  # I132: 1351.4259 -> 1917.7667
$ws.Cells.Item(132, 11).Value = 5753.300099999999  # K132: 4054.2777 -> 5753.300099999999
$ws.Cells.Item(132, 13).Value = -3223.300099999999  # M132: -1524.2777 -> -3223.300099999999
